# Trim the "Rendimento " prefix from the two column headers on the
# "Séries" sheet (B1/C1), and move the active-cell selection from C1 to B1
# to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "efetivo real"
$ws.Range("C1").Value = "habitual real"

$ws.Activate() | Out-Null
$ws.Range("B1").Select() | Out-Null
